# "Added the reporting piece"
#
# The "Result" header is renamed to "Results" on the "Test Cases" and
# "GoogleSearch" sheets, and a "PASS" status value is recorded under the
# renamed header on both sheets (these two sheets previously left that
# column empty).

$wb = $excel.ActiveWorkbook

# --- "Test Cases" sheet ---------------------------------------------------
$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsTestCases.Range("E1").Value = "Results"
$wsTestCases.Range("E2").Value = "PASS"
$wsTestCases.Columns.Item(5).AutoFit()
[void]$wsTestCases.Range("E2").Select()

# --- "GoogleSearch" sheet --------------------------------------------------
$wsGoogleSearch = $wb.Worksheets.Item("GoogleSearch")
$wsGoogleSearch.Range("D2").Value = "PASS"
$wsGoogleSearch.Range("D1").Value = "Results"
$wsGoogleSearch.Columns.Item(4).AutoFit()
[void]$wsGoogleSearch.Range("D1").Select()

$wb.Save()
